$ws = $excel.ActiveWorkbook.ActiveSheet

# Rename the "Rolling_Avg_Misaligned_1M" label to "_6M" (header C1 and row-label A3)
$ws.Range("C1").Value = "Rolling_Avg_Misaligned_6M"
$ws.Range("A3").Value = "Rolling_Avg_Misaligned_6M"

# Updated correlation values (recomputed matrix -- new 6M column/row plus refreshed float precision)
$ws.Range("C2").Value = 0.02435689520289818
$ws.Range("D2").Value = 0.04226782259869652
$ws.Range("H2").Value = -0.00144532414223878
$ws.Range("K2").Value = 0.05668574567102287
$ws.Range("N2").Value = -0.08887914460199312
$ws.Range("B3").Value = 0.02435689520289818
$ws.Range("D3").Value = -0.06928357492136708
$ws.Range("E3").Value = 0.1618402499033878
$ws.Range("F3").Value = -0.01858395782870342
$ws.Range("G3").Value = 0.003482931764957531
$ws.Range("H3").Value = 0.001884881086372575
$ws.Range("I3").Value = 0.06744470724421375
$ws.Range("J3").Value = 0.04841034537771964
$ws.Range("K3").Value = 0.05565677493390313
$ws.Range("L3").Value = 0.01217573226643082
$ws.Range("M3").Value = 0.04606925186632603
$ws.Range("N3").Value = 0.04579195693789806
$ws.Range("B4").Value = 0.04226782259869652
$ws.Range("C4").Value = -0.06928357492136708
$ws.Range("E4").Value = 0.06149912334035905
$ws.Range("F4").Value = 0.02413345473059602
$ws.Range("H4").Value = -0.002372091436832229
$ws.Range("I4").Value = 0.0549593996644491
$ws.Range("K4").Value = 0.04014502415598669
$ws.Range("N4").Value = -0.01893148487609338
$ws.Range("C5").Value = 0.1618402499033878
$ws.Range("D5").Value = 0.06149912334035905
$ws.Range("F5").Value = 0.01924531793298757
$ws.Range("G5").Value = -0.01081426358502754
$ws.Range("H5").Value = -0.0007567622585659154
$ws.Range("J5").Value = -0.02065820748257093
$ws.Range("K5").Value = 0.07480425044441923
$ws.Range("M5").Value = -0.004832613374206191
$ws.Range("C6").Value = -0.01858395782870342
$ws.Range("D6").Value = 0.02413345473059602
$ws.Range("E6").Value = 0.01924531793298757
$ws.Range("G6").Value = 0.08181963744959719
$ws.Range("H6").Value = -0.0009238941443531934
$ws.Range("J6").Value = -0.003482285323228752
$ws.Range("L6").Value = -0.09351902174369986
$ws.Range("C7").Value = 0.003482931764957531
$ws.Range("E7").Value = -0.01081426358502754
$ws.Range("F7").Value = 0.08181963744959719
$ws.Range("H7").Value = 0.002383877351227499
$ws.Range("J7").Value = -0.05987449686194377
$ws.Range("K7").Value = 0.3798642189897706
$ws.Range("M7").Value = 0.03640899503708583
$ws.Range("B8").Value = -0.00144532414223878
$ws.Range("C8").Value = 0.001884881086372575
$ws.Range("D8").Value = -0.002372091436832229
$ws.Range("E8").Value = -0.0007567622585659154
$ws.Range("F8").Value = -0.0009238941443531934
$ws.Range("G8").Value = 0.002383877351227499
$ws.Range("I8").Value = 0.000009611274414439937
$ws.Range("J8").Value = 0.001063558458272517
$ws.Range("K8").Value = -0.002367311497960982
$ws.Range("L8").Value = 0.0006593989500711931
$ws.Range("M8").Value = -0.003756267431308731
$ws.Range("N8").Value = 0.0009322909582774932
$ws.Range("C9").Value = 0.06744470724421375
$ws.Range("D9").Value = 0.0549593996644491
$ws.Range("H9").Value = 0.000009611274414439937
$ws.Range("J9").Value = 0.01978848351359107
$ws.Range("K9").Value = 0.1854957068454134
$ws.Range("L9").Value = -0.06405003630667885
$ws.Range("M9").Value = 0.05400895170944196
$ws.Range("N9").Value = -0.0264255798556076
$ws.Range("C10").Value = 0.04841034537771964
$ws.Range("E10").Value = -0.02065820748257093
$ws.Range("F10").Value = -0.003482285323228752
$ws.Range("G10").Value = -0.05987449686194377
$ws.Range("H10").Value = 0.001063558458272517
$ws.Range("I10").Value = 0.01978848351359107
$ws.Range("N10").Value = 0.01892521710331613
$ws.Range("B11").Value = 0.05668574567102287
$ws.Range("C11").Value = 0.05565677493390313
$ws.Range("D11").Value = 0.04014502415598669
$ws.Range("E11").Value = 0.07480425044441923
$ws.Range("G11").Value = 0.3798642189897706
$ws.Range("H11").Value = -0.002367311497960982
$ws.Range("I11").Value = 0.1854957068454134
$ws.Range("M11").Value = 0.2987826524325734
$ws.Range("N11").Value = -0.1653177463296434
$ws.Range("C12").Value = 0.01217573226643082
$ws.Range("F12").Value = -0.09351902174369986
$ws.Range("H12").Value = 0.0006593989500711931
$ws.Range("I12").Value = -0.06405003630667885
$ws.Range("C13").Value = 0.04606925186632603
$ws.Range("E13").Value = -0.004832613374206191
$ws.Range("G13").Value = 0.03640899503708583
$ws.Range("H13").Value = -0.003756267431308731
$ws.Range("I13").Value = 0.05400895170944196
$ws.Range("K13").Value = 0.2987826524325734
$ws.Range("B14").Value = -0.08887914460199312
$ws.Range("C14").Value = 0.04579195693789806
$ws.Range("D14").Value = -0.01893148487609338
$ws.Range("H14").Value = 0.0009322909582774932
$ws.Range("I14").Value = -0.0264255798556076
$ws.Range("J14").Value = 0.01892521710331613
$ws.Range("K14").Value = -0.1653177463296434
